$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: component/header block becomes a multi-line placeholder block ---
$ws.Range("B2").Value = "Component: MultiFunctionalTool`nMFP: Any`nBuild: `nDate: `nTarget: "

# --- Rows 6-21: wipe out the generated test-case rows, leave a "Not Executed" marker in G ---
for ($r = 6; $r -le 21; $r++) {
    $ws.Range("B$r`:F$r").ClearContents()
    $ws.Range("G$r").Value = "Not Executed"
    $ws.Range("H$r").ClearContents()
}

# --- Row 22: also wiped (row kept, since its cells carry explicit styles) ---
$ws.Range("B22:H22").ClearContents()

# --- Row 23: wiped entirely (cells had no explicit style, so the row disappears) ---
$ws.Range("B23:H23").ClearContents()

# --- Row 24: becomes the "Test Summary" label row ---
$ws.Range("B24").Value = "Test Summary"
$ws.Range("C24:H24").ClearContents()

# --- Row 25: becomes the "Test Case Count:" label row ---
$ws.Range("D25").Value = "Test Case Count:"
$ws.Range("B25:C25").ClearContents()
$ws.Range("E25:H25").ClearContents()

# --- Row 26: becomes the "Total Test Cases" formula row ---
$ws.Range("D26").Value = "Total Test Cases"
$ws.Range("B26:C26").ClearContents()
$ws.Range("E26").Formula = "=COUNTA(B6:B6)"
$ws.Range("F26:H26").ClearContents()

# --- Row 27: becomes the "Total Test Case Passed" formula row ---
$ws.Range("D27").Value = "Total Test Case Passed"
$ws.Range("B27:C27").ClearContents()
$ws.Range("E27").Formula = "=COUNTIF(G6:G6,""PASS"")"
$ws.Range("F27:H27").ClearContents()
